$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-gender")

# Update is_active (column D) to FALSE for the "Others" (OTH) rows
$ws.Range("D4").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("D10").Value = $false

# Update the active selection to D12
$ws.Range("D12").Select()
